# Case with 380 kV: refreshed simulation result values in pl_mw.xlsx
# (sheet1, rows 2-25, columns B/D/E/F/G/H/J/M/N). Columns A, C, I, K, L, O
# are left untouched since their values did not change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    'B2' = 0.1424644886134701; 'D2' = 0.1139241056133855; 'E2' = 0.1040233722974415; 'F2' = 2.241138824841414; 'G2' = 1.776117867831744; 'H2' = 1.40473911688116; 'J2' = 0.1249938732257547; 'M2' = 1.206681053132939; 'N2' = 1.459012336282058
    'B3' = 0.1329234947758806; 'D3' = 0.1119549622607394; 'E3' = 0.1045922858359685; 'F3' = 2.184739009880289; 'G3' = 1.698676187230404; 'H3' = 1.376547655409496; 'J3' = 0.1270065260697453; 'M3' = 1.096913904105591; 'N3' = 1.423636744251979
    'B4' = 0.1271321576970905; 'D4' = 0.1107307092836649; 'E4' = 0.10496619120826; 'F4' = 2.151842812373971; 'G4' = 1.652587754564991; 'H4' = 1.360257592255891; 'J4' = 0.128309043113616; 'M4' = 1.029692168599212; 'N4' = 1.402506520986677
    'B5' = 0.1247891573038231; 'D5' = 0.1102279396534129; 'E5' = 0.105124761798626; 'F5' = 2.138869424792119; 'G5' = 1.634169209597445; 'H5' = 1.353873714725268; 'J5' = 0.1288565977352847; 'M5' = 1.002343018362609; 'N5' = 1.394045718781229
    'B6' = 0.1244011375526952; 'D6' = 0.110144219361068; 'E6' = 0.1051514673760154; 'F6' = 2.136741189375172; 'G6' = 1.631132609032875; 'H6' = 1.352828989243363; 'J6' = 0.1289485311254523; 'M6' = 0.9978043882831003; 'N6' = 1.39264990751407
    'B7' = 0.1271004900081891; 'D7' = 0.1107239445274573; 'E7' = 0.1049683046148342; 'F7' = 2.151666104264322; 'G7' = 1.652337892338466; 'H7' = 1.360170469296719; 'J7' = 0.1283163597496406; 'M7' = 1.029323149106204; 'N7' = 1.402391806635706
    'B8' = 0.1391609919028838; 'D8' = 0.113248263070254; 'E8' = 0.1042144406700745; 'F8' = 2.221330035507549; 'G8' = 1.749109958569193; 'H8' = 1.394805896804314; 'J8' = 0.1256739757904306; 'M8' = 1.168796857086448; 'N8' = 1.446693110347894
    'B9' = 0.1633346837021747; 'D9' = 0.1180806803733816; 'E9' = 0.1029304327905258; 'F9' = 2.371879878529853; 'G9' = 1.95069064692953; 'H9' = 1.470909624759202; 'J9' = 0.1210223732031794; 'M9' = 1.443712307399011; 'N9' = 1.538196776770633
    'B10' = 0.1814058279748423; 'D10' = 0.1215634100430734; 'E10' = 0.1021044401319173; 'F10' = 2.491258968324189; 'G10' = 2.106313795450774; 'H10' = 1.531951058271147; 'J10' = 0.1179286587131241; 'M10' = 1.646593425051563; 'N10' = 1.608180019640741
    'B11' = 0.1896927262001356; 'D11' = 0.1231340313962477; 'E11' = 0.1017539269845367; 'F11' = 2.547532837944487; 'G11' = 2.178815780071091; 'H11' = 1.560865031734465; 'J11' = 0.1165917047033771; 'M11' = 1.739096157013819; 'N11' = 1.640602413352838
    'B12' = 0.1928401076212367; 'D12' = 0.1237268868753958; 'E12' = 0.1016248071457744; 'F12' = 2.569130019111299; 'G12' = 2.206521643586882; 'H12' = 1.571981223268835; 'J12' = 0.1160955738491523; 'M12' = 1.774155401764489; 'N12' = 1.652963051266852
    'B13' = 0.192161851909006; 'D13' = 0.1235992885245594; 'E13' = 0.1016524550426094; 'F13' = 2.564465820511117; 'G13' = 2.200543444410016; 'H13' = 1.569579683439201; 'J13' = 0.1162019729612407; 'M13' = 1.766603408732152; 'N13' = 1.650297297859794
    'B14' = 0.1899514777369689; 'D14' = 0.1231828437763767; 'E14' = 0.1017432319267391; 'F14' = 2.549303863081803; 'G14' = 2.181090098351547; 'H14' = 1.561776205193269; 'J14' = 0.1165506842771746; 'M14' = 1.741979890738946; 'N14' = 1.641617675234613
    'B15' = 0.1885987662387834; 'D15' = 0.1229275132079124; 'E15' = 0.1017993052622961; 'F15' = 2.540054301775484; 'G15' = 2.169207200144172; 'H15' = 1.557018179973909; 'J15' = 0.1167656018242571; 'M15' = 1.726901246966349; 'N15' = 1.636311921764928
    'B16' = 0.1808655744775649; 'D16' = 0.121460496197642; 'E16' = 0.1021278533403533; 'F16' = 2.487621319639999; 'G16' = 2.101610444591188; 'H16' = 1.530084717787418; 'J16' = 0.1180174501690647; 'M16' = 1.640552410229361; 'N16' = 1.606072826605953
    'B17' = 0.1761383268766821; 'D17' = 0.1205570675008971; 'E17' = 0.1023358582270593; 'F17' = 2.455962561668969; 'G17' = 2.060583280324323; 'H17' = 1.513857041121952; 'J17' = 0.1188034634424162; 'M17' = 1.5876343265939; 'N17' = 1.587671390390909
    'B18' = 0.1734255869915842; 'D18' = 0.1200361467608175; 'E18' = 0.1024578733885004; 'F18' = 2.437938216636326; 'G18' = 2.03714600530455; 'H18' = 1.504631028758013; 'J18' = 0.1192621817798161; 'M18' = 1.557217103894558; 'N18' = 1.577142671766069
    'B19' = 0.1725081794553063; 'D19' = 0.1198595480454259; 'E19' = 0.1024995942638274; 'F19' = 2.431867109206451; 'G19' = 2.029237946925718; 'H19' = 1.501525690367259; 'J19' = 0.1194186329818705; 'M19' = 1.546921763169635; 'N19' = 1.573587368961626
    'B20' = 0.176640905362845; 'D20' = 0.1206533723332655; 'E20' = 0.1023134699411692; 'F20' = 2.459313517861744; 'G20' = 2.064934046472644; 'H20' = 1.515573339977749; 'J20' = 0.1187191052382173; 'M20' = 1.593265490904329; 'N20' = 1.589624541345501
    'B21' = 0.1906004671568269; 'D21' = 0.1233052148729499; 'E21' = 0.1017164706645088; 'F21' = 2.553749458615073; 'G21' = 2.1867971646613; 'H21' = 1.56406372224842; 'J21' = 0.116447983808587; 'M21' = 1.749211586827045; 'N21' = 1.644164849601196
    'B22' = 0.1997780320469218; 'D22' = 0.1250272755548139; 'E22' = 0.1013473421198157; 'F22' = 2.617146671080548; 'G22' = 2.267906926976821; 'H22' = 1.596730069809041; 'J22' = 0.1150228159002671; 'M22' = 1.85130905254502; 'N22' = 1.680293146946696
    'B23' = 0.1948749053101579; 'D23' = 0.1241091709877722; 'E23' = 0.1015424329279995; 'F23' = 2.583155340605742; 'G23' = 2.224481282848899; 'H23' = 1.579205427962279; 'J23' = 0.1157780359230869; 'M23' = 1.796801405984155; 'N23' = 1.66096704545285
    'B24' = 0.1764136739845128; 'D24' = 0.1206098376821387; 'E24' = 0.1023235841138499; 'F24' = 2.457798000538105; 'G24' = 2.062966598939255; 'H24' = 1.51479707907589; 'J24' = 0.1187572223139339; 'M24' = 1.590719622562617; 'N24' = 1.588741364359691
    'B25' = 0.1567399045443665; 'D25' = 0.1167855665534461; 'E25' = 0.103257100983325; 'F25' = 2.329632003430135; 'G25' = 1.894860365773866; 'H25' = 1.449431212419285; 'J25' = 0.122223937791988; 'M25' = 1.369185445620488; 'N25' = 1.512954008150672
}

foreach ($cellAddr in $newValues.Keys) {
    $ws.Range($cellAddr).Value = [double]$newValues[$cellAddr]
}
